# Regenerate save_data to use K instead of Strike#: update column G (K) values
# for rows 2-29 on the active sheet, matching the recalculated stat values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 1
    4  = 3
    5  = 1
    6  = 4
    7  = 5
    8  = 3
    9  = 6
    10 = 1
    11 = 2
    12 = 2
    13 = 1
    14 = 2
    15 = 4
    16 = 4
    17 = 1
    18 = 3
    19 = 2
    20 = 4
    21 = 7
    22 = 0
    23 = 4
    24 = 5
    25 = 6
    26 = 3
    27 = 3
    28 = 2
    29 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

Write-Host "Updated column G (K) values for rows 2-29"
